$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.169.71"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.325.18"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.15"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.63"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.319.32"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.178"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.576"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.58"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000266"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.859.57"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.49"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "605.70"
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.286.76"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.93"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.329.49"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.09"
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.04"
$ws.Range("E24").Value = "  -4.64%  "
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.01"
$ws.Range("E26").Value = "  +3.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.73"
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.54"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "31.03"
$ws.Range("E29").Value = "  +3.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.53"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.57"
$ws.Range("E31").Value = "  +4.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.73"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "556.86"
$ws.Range("E33").Value = "  +4.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.92"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.828.60"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.104"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.28"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.129"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.11"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.17"
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0686"
$ws.Range("E42").Value = "  -5.66%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.36"
$ws.Range("E44").Value = "  +4.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.00"
$ws.Range("E47").Value = "  -7.39%  "
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.41"
$ws.Range("E51").Value = "  +6.01%  "
